# 学習記録.xlsx - add new training-log entry row (2023-12-19, UniPC / PPO unit-count
# experiment) and a new "追記" remarks column, per commit "changed unit num of PPO".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: apply the CJK body font to a whole cell (matches the sheet's
# "Noto Sans CJK SC" 10pt style used throughout column headers / remarks).
# ---------------------------------------------------------------------------
function Set-CjkFont($rng) {
    $rng.Font.Name = "Noto Sans CJK SC"
    $rng.Font.Size = 10
    $rng.Font.ColorIndex = -4105
}

function Set-LatinFont($rng) {
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 10
    $rng.Font.ColorIndex = -4105
    $rng.Font.Charset = 1
}

# ---------------------------------------------------------------------------
# Header row: new "追記" column G, and F1 ("備考") gets the other CJK font
# slot (the one already used by the rest of the header row / E1 etc).
# ---------------------------------------------------------------------------
Set-CjkFont($ws.Range("F1"))

$ws.Range("G1").Value = "追記"
Set-CjkFont($ws.Range("G1"))

# F2 ("ドアヒンジの...") moves to the same font slot as F1.
Set-CjkFont($ws.Range("F2"))

# ---------------------------------------------------------------------------
# New row 3: 12/19 UniPC entry about increasing the PPO MLP unit count.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 45279
$ws.Range("A3").NumberFormat = "mm/dd/yy"

$ws.Range("B3").Value = "UniPC"

$ws.Range("C3").Value = "同上"
Set-CjkFont($ws.Range("C3"))

# D3: mixed Latin/CJK rich text run-by-run.
$ws.Range("D3").Value = "mlpのunitを増やしてみた5層，2048→64"
Set-LatinFont($ws.Range("D3").Characters(1, 3))      # "mlp"
Set-CjkFont($ws.Range("D3").Characters(4, 1))        # "の"
Set-LatinFont($ws.Range("D3").Characters(5, 4))      # "unit"
Set-CjkFont($ws.Range("D3").Characters(9, 7))        # "を増やしてみた"
Set-LatinFont($ws.Range("D3").Characters(16, 1))     # "5"
Set-CjkFont($ws.Range("D3").Characters(17, 2))       # "層，"
Set-LatinFont($ws.Range("D3").Characters(19, 7))     # "2048→64"

$ws.Range("E3").Value = "✕ 要検討"
Set-CjkFont($ws.Range("E3"))

# F3: mixed Latin/CJK rich text run-by-run.
$ws.Range("F3").Value = "unit数増やしすぎても逆効果かも PPOだと報酬の変化方向に動くはずなのでスパースだと無理な気がしてきた"
Set-LatinFont($ws.Range("F3").Characters(1, 4))      # "unit"
Set-CjkFont($ws.Range("F3").Characters(5, 14))       # "数増やしすぎても逆効果かも "
Set-LatinFont($ws.Range("F3").Characters(19, 3))     # "PPO"
Set-CjkFont($ws.Range("F3").Characters(22, 32))      # "だと報酬の変化方向に動くはずなのでスパースだと無理な気がしてきた"

$ws.Range("G3").Value = "ゴミなので特にデータ出力せず"
Set-CjkFont($ws.Range("G3"))

# ---------------------------------------------------------------------------
# Column widths / sheet view cosmetics.
# ---------------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 93.77
$ws.Columns("G").ColumnWidth = 18.82

$excel.ActiveWindow.Zoom = 140
[void]$ws.Range("G4").Select()

Write-Host "edit complete"
